$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 182; everything from the old row 182
# onward shifts down by one (old row 284 -> new row 285).
$ws.Rows.Item(182).Insert()

# Populate the newly inserted row 182 with the new record.
$ws.Range("A182").Value = 3
$ws.Range("B182").Value = "Femacal de La Calera"
$ws.Range("C182").Value = "Coquimbo"
$ws.Range("D182").Value = 44529
$ws.Range("E182").Value = 5
$ws.Range("F182").Value = 100112028
$ws.Range("G182").Value = "Sandia"
$ws.Range("H182").Value = "Sin especificar"
$ws.Range("I182").Value = "Primera"
$ws.Range("J182").Value = 270
$ws.Range("K182").Value = 700
$ws.Range("L182").Value = 750
$ws.Range("M182").Value = 720
$ws.Range("N182").Value = "$/kilo (volumen en unidades)"
$ws.Range("O182").Value = "Perú"
$ws.Range("P182").Value = 720
$ws.Range("Q182").Value = 1
$ws.Range("R182").Value = "Hortaliza"
